$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2 and 3 with new capital structure values
foreach ($r in 2,3) {
    $ws.Range("K$r").Value = -2.06
    $ws.Range("U$r").Value = 0.364
    $ws.Range("V$r").Value = 0.01016759776536313
    $ws.Range("W$r").Value = -0.05309278350515464
    $ws.Range("X$r").Value = 0.0842299367121281
    $ws.Range("Y$r").Value = -0.1373227202172828
    $ws.Range("AA$r").Value = -0.03844100373731981
    $ws.Range("AB$r").Value = 0.0842299367121281
    $ws.Range("AC$r").Value = -0.1226709404494479
    $ws.Range("AG$r").Value = -0.364
    $ws.Range("AJ$r").Value = -0.01027203973360424
    $ws.Range("AK$r").Value = -0.009620467279839305
    $ws.Range("AM$r").Value = -0.004
    $ws.Range("AQ$r").Value = 360
}
